# Scheduled market-data refresh: updates the currentAveragePrice* /
# LevePrice* / LeveProfit* columns (H:N) on every Leve worksheet with
# freshly polled prices. Columns A:G (leve metadata) are untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 728.72
$ws.Range("J17").Value = 728.72
$ws.Range("L17").Value = 2186.16
$ws.Range("N17").Value = -2522.16

$ws.Range("H19").Value = 5445.75
$ws.Range("J19").Value = 6892
$ws.Range("L19").Value = 6892
$ws.Range("N19").Value = -7242

$ws.Range("H33").Value = 450.1579
$ws.Range("I33").Value = 265.84616
$ws.Range("J33").Value = 849.5
$ws.Range("K33").Value = 265.84616
$ws.Range("L33").Value = 849.5
$ws.Range("M33").Value = -36.84616
$ws.Range("N33").Value = -1307.5

$ws.Range("H86").Value = 1795
$ws.Range("I86").Value = 1519.8572
$ws.Range("K86").Value = 1519.8572
$ws.Range("M86").Value = -396.8571999999999

$ws.Range("H89").Value = 1795
$ws.Range("I89").Value = 1519.8572
$ws.Range("K89").Value = 7599.286
$ws.Range("M89").Value = -1983.286

$ws.Range("H100").Value = 1376.4615
$ws.Range("I100").Value = 1254
$ws.Range("K100").Value = 1254
$ws.Range("M100").Value = -713

$ws.Range("H103").Value = 1473.3572
$ws.Range("I103").Value = 2122.1667
$ws.Range("K103").Value = 6366.500100000001
$ws.Range("M103").Value = -5780.500100000001

$ws.Range("H113").Value = 4572.222
$ws.Range("I113").Value = 2351.25
$ws.Range("J113").Value = 6349
$ws.Range("K113").Value = 2351.25
$ws.Range("L113").Value = 6349
$ws.Range("M113").Value = 902.75
$ws.Range("N113").Value = -12857

$ws.Range("H116").Value = 13033.05
$ws.Range("I116").Value = 13033.05
$ws.Range("K116").Value = 13033.05
$ws.Range("M116").Value = -9591.049999999999

$ws.Range("H135").Value = 1233.4117
$ws.Range("I135").Value = 597.86664
$ws.Range("K135").Value = 5380.79976
$ws.Range("M135").Value = -2845.79976

$ws.Range("H137").Value = 15122.913
$ws.Range("I137").Value = 19319.53
$ws.Range("K137").Value = 57958.59
$ws.Range("M137").Value = -55408.59

$ws.Range("H138").Value = 27191.365
$ws.Range("I138").Value = 2113.7
$ws.Range("J138").Value = 51074.855
$ws.Range("K138").Value = 6341.099999999999
$ws.Range("L138").Value = 153224.565
$ws.Range("M138").Value = -1201.099999999999
$ws.Range("N138").Value = -163504.565

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 660
$ws.Range("I2").Value = 660
$ws.Range("K2").Value = 660
$ws.Range("M2").Value = -547

$ws.Range("H32").Value = 34763.406
$ws.Range("I32").Value = 34763.406
$ws.Range("K32").Value = 34763.406
$ws.Range("M32").Value = -34476.406

$ws.Range("H61").Value = 4749.2856
$ws.Range("I61").Value = 915.875
$ws.Range("K61").Value = 915.875
$ws.Range("M61").Value = -703.875

$ws.Range("H74").Value = 322006.38
$ws.Range("I74").Value = 400714.72
$ws.Range("J74").Value = 26850
$ws.Range("K74").Value = 400714.72
$ws.Range("L74").Value = 26850
$ws.Range("M74").Value = -399840.72
$ws.Range("N74").Value = -28598

$ws.Range("H77").Value = 322006.38
$ws.Range("I77").Value = 400714.72
$ws.Range("J77").Value = 26850
$ws.Range("K77").Value = 2003573.6
$ws.Range("L77").Value = 134250
$ws.Range("M77").Value = -1999205.6
$ws.Range("N77").Value = -142986

$ws.Range("H96").Value = 60000
$ws.Range("J96").Value = 60000
$ws.Range("L96").Value = 60000
$ws.Range("N96").Value = -65492

$ws.Range("H116").Value = 660
$ws.Range("I116").Value = 660
$ws.Range("K116").Value = 660
$ws.Range("M116").Value = 1634

$ws.Range("H132").Value = 1445.3733
$ws.Range("I132").Value = 1063.1428
$ws.Range("K132").Value = 3189.4284
$ws.Range("M132").Value = -659.4284000000002

$ws.Range("H136").Value = 4749.2856
$ws.Range("I136").Value = 915.875
$ws.Range("K136").Value = 2747.625
$ws.Range("M136").Value = -197.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 660
$ws.Range("I3").Value = 660
$ws.Range("K3").Value = 660
$ws.Range("M3").Value = -546

$ws.Range("H86").Value = 4067.5715
$ws.Range("I86").Value = 1699.6666
$ws.Range("K86").Value = 1699.6666
$ws.Range("M86").Value = -576.6666

$ws.Range("H89").Value = 4067.5715
$ws.Range("I89").Value = 1699.6666
$ws.Range("K89").Value = 8498.333000000001
$ws.Range("M89").Value = -2882.333000000001

$ws.Range("H99").Value = 2084.9333
$ws.Range("I99").Value = 1948.5834
$ws.Range("K99").Value = 1948.5834
$ws.Range("M99").Value = -450.5834

$ws.Range("H107").Value = 6558.273
$ws.Range("I107").Value = 6547.5
$ws.Range("K107").Value = 6547.5
$ws.Range("M107").Value = -4627.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4349977.5
$ws.Range("I31").Value = 6250998
$ws.Range("K31").Value = 6250998
$ws.Range("M31").Value = -6250703

$ws.Range("H34").Value = 4349977.5
$ws.Range("I34").Value = 6250998
$ws.Range("K34").Value = 6250998
$ws.Range("M34").Value = -6250796

$ws.Range("H105").Value = 56208.855
$ws.Range("I105").Value = 64743.668
$ws.Range("K105").Value = 64743.668
$ws.Range("M105").Value = -62996.668

$ws.Range("H141").Value = 66803.836
$ws.Range("J141").Value = 70164.60000000001
$ws.Range("L141").Value = 70164.60000000001
$ws.Range("N141").Value = -80524.60000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 204.1
$ws.Range("J12").Value = 125.125
$ws.Range("L12").Value = 375.375
$ws.Range("N12").Value = -721.375

$ws.Range("H37").Value = 42539
$ws.Range("J37").Value = 42539
$ws.Range("L37").Value = 127617
$ws.Range("N37").Value = -127841

$ws.Range("H131").Value = 477516.56
$ws.Range("J131").Value = 3666
$ws.Range("L131").Value = 10998
$ws.Range("N131").Value = -21078

$ws.Range("H132").Value = 1229.7273
$ws.Range("I132").Value = 900
$ws.Range("J132").Value = 1353.375
$ws.Range("K132").Value = 8100
$ws.Range("L132").Value = 12180.375
$ws.Range("M132").Value = -5570
$ws.Range("N132").Value = -17240.375

$ws.Range("H134").Value = 684.3570999999999
$ws.Range("I134").Value = 349.84616
$ws.Range("K134").Value = 1049.53848
$ws.Range("M134").Value = 4020.46152

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 2933.3333
$ws.Range("I9").Value = 1900
$ws.Range("J9").Value = 5000
$ws.Range("K9").Value = 1900
$ws.Range("L9").Value = 5000
$ws.Range("M9").Value = -1730
$ws.Range("N9").Value = -5340

$ws.Range("H122").Value = 3632.9443
$ws.Range("I122").Value = 3520.9167
$ws.Range("J122").Value = 3857
$ws.Range("K122").Value = 10562.7501
$ws.Range("L122").Value = 11571
$ws.Range("M122").Value = -8112.750100000001
$ws.Range("N122").Value = -16471

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2261.3572
$ws.Range("I7").Value = 2288.111
$ws.Range("K7").Value = 2288.111
$ws.Range("M7").Value = -2176.111

$ws.Range("H126").Value = 2261.3572
$ws.Range("I126").Value = 2288.111
$ws.Range("K126").Value = 6864.333
$ws.Range("M126").Value = -4394.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6887.1113
$ws.Range("I81").Value = 9099.799999999999
$ws.Range("J81").Value = 4121.25
$ws.Range("K81").Value = 18199.6
$ws.Range("L81").Value = 8242.5
$ws.Range("M81").Value = -17138.6
$ws.Range("N81").Value = -10364.5

$ws.Range("H84").Value = 6887.1113
$ws.Range("I84").Value = 9099.799999999999
$ws.Range("J84").Value = 4121.25
$ws.Range("K84").Value = 90998
$ws.Range("L84").Value = 41212.5
$ws.Range("M84").Value = -85694
$ws.Range("N84").Value = -51820.5

$ws.Range("H126").Value = 195433.34
$ws.Range("I126").Value = 2842.1667
$ws.Range("K126").Value = 8526.500100000001
$ws.Range("M126").Value = -6056.500100000001

$ws.Range("H132").Value = 1682.7407
$ws.Range("I132").Value = 1221.8
$ws.Range("K132").Value = 3665.4
$ws.Range("M132").Value = -1135.4
